$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.074.83"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "1.826.86"
$ws.Range("E3").Value = "  +0.43%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.010"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.77%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "312.71"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.80%  "
$ws.Range("E6").Value = "  +0.71%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4697"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.22%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3647"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.47%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07387"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +0.46%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.8810"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.00%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "20.30"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").Value = "1.933.91"
$ws.Range("E12").Value = "  +6.42%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.07313"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +2.91%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "93.33"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +2.03%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "5.365"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.84%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "6.529"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.31%  "
$ws.Range("D17").Value = "1.008"
$ws.Range("E17").Value = "  +0.48%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.000008709"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.80%  "
$ws.Range("D20").Value = "27.594.96"
$ws.Range("E20").Value = "  +2.26%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "14.63"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.32%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.241"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.85%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "10.63"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("D24").Value = "2.096.26"
$ws.Range("E24").Value = "  +2.60%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "1.884"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.48%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "151.73"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.42%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "18.50"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.80%  "
$ws.Range("D28").Value = "2.137"
$ws.Range("E28").Value = "  -0.74%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "5.190"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -1.11%  "
$ws.Range("D30").Value = "116.47"
$ws.Range("E30").Value = "  -0.38%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.08941"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.63%  "
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("D33").Value = "0.7409"
$ws.Range("E33").Value = "  -2.31%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "4.511"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.05%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "2.946"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.64%  "
$ws.Range("E36").Value = "  +0.79%  "
$ws.Range("D37").Value = "1.089"
$ws.Range("E37").Value = "  -0.67%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.05298"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("E39").Value = "  -0.10%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.405"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +1.59%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "2.928"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -1.29%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "7.214"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.45%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.5245"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.81%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.1642"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.69%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "8.390"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.61%  "
$ws.Range("D46").Value = "0.4874"
$ws.Range("E46").Value = "  +0.02%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "10.38"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.62%  "
$ws.Range("E48").Value = "  +0.77%  "
$ws.Range("D49").Value = "104.34"
$ws.Range("E49").Value = "  +0.90%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.648"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.73%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.06284"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.22%  "
